# Work on the stop signal task
#
# - Re-crop / reposition the top "person at computer" photo (Picture 2):
#   the source-rectangle top crop grows from 9.142% to 15.28%, the photo
#   shifts up slightly and gets a bit shorter.
# - Shrink the description textbox (TextBox 11) to match the photo's new
#   (smaller) footprint.
# - Rewrite the first descriptive paragraph's wording.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: "Picture 2" (top photo) ---------------------------------
$pic = $s.Shapes.Item(1)
$pic.PictureFormat.CropTop = 283.635
$pic.Top = -8.0
$pic.Height = 457.4857480314961
# width/left are unchanged (540 / 0), left untouched intentionally

# --- Shape 3: "TextBox 11" (description box) ---------------------------
$textbox = $s.Shapes.Item(3)
$textbox.Height = 230.22653543307086

$para = $textbox.TextFrame.TextRange.Paragraphs(4)
$run = $para.Runs(1)
$run.Text = "The Stop-Signal Task (SST) is a psychological experiment designed to study response inhibition, a key aspect of cognitive control. In this task, participants are asked to respond quickly to a `"go`" signal, by pressing a button when they see an arrow pointing left or right. However, on some trials, a `"stop`" signal (a red round circle underneath the arrow) appears shortly after the go signal. In these trials the participant must withhold their response."
